$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 2.2
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 6
$ws.Range("AC5").Value = 7
$ws.Range("AE5").Value = 23
$ws.Range("AH5").Value = 34
$ws.Range("AJ5").Value = 81
$ws.Range("AN5").Value = 3.25
$ws.Range("AS5").Value = 251
$ws.Range("AX5").Value = 41

# Row 7 updates
$ws.Range("G7").Value = 2.25
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 3.1
$ws.Range("X7").Value = 10
$ws.Range("AK7").Value = 29
$ws.Range("AX7").Value = 19
$ws.Range("BB7").Value = 251
